$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Test-Payload" for row 20: temperature JSON gains an extra "unexpected" field
$newTestPayload = @"
{
     "temperature": 120.5,
     "unit": "Celsius",
     "time": "2023-07-12T16:21:53.389+02:00",
     "externalId": "berlin_01",
     "unexpected": 17.5
}
"@

# New "Expected Result" for row 20: extra paragraph explaining the CREATE_IF_MISSING repair strategy
$newExpectedResult = @"
A measasurement should be created for the device berlin_01.
The fragment "c8y_Fragment_to_remove" is not included in the created measurement, as the repair strategy is "REMOVE_IF_NULL".
In addition the reapar strategy "CREATE_IF_MISSING" is used. Thjsi is required to map the node "unexpected" to the target fragment "c8y_Unexpected". This is created, due to the used reapir strategy.
"@

# Row 20 previously held the payload/description pair in E20/F20; the diff swaps which
# column carries which payload: the (now extended) Test-Payload moves to E20 while the
# Target-Payload (unchanged "source" json, previously in F20) stays content-identical in F20.
$ws.Range("E20").Value = $newTestPayload
$ws.Range("F20").Value = $ws.Range("F20").Value2
$ws.Range("G20").Value = $newExpectedResult

# Row 20 grows taller to fit the additional text that was added.
$ws.Rows.Item(20).RowHeight = 252

# Rows 2 and 3 are re-clamped from 409.6 to Excel's real max row height of 409.5.
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5

# Reflect the updated view/selection state pointing at the edited row.
$ws.Application.Goto($ws.Range("B18"), $true)
$ws.Range("F20").Select()
